$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = -1
$ws.Cells.Item(3,1).Value = -0.00038742819433336596
$ws.Cells.Item(3,2).Value = -0.06479935270847112
$ws.Cells.Item(3,3).Value = 0.9978982381926312
$ws.Cells.Item(4,1).Value = 0.0007653948946984158
$ws.Cells.Item(4,2).Value = 0.12801546965896227
$ws.Cells.Item(4,3).Value = -0.9917718758356936
$ws.Cells.Item(5,1).Value = -0.0011555889232988086
$ws.Cells.Item(5,2).Value = -0.19327813143520226
$ws.Cells.Item(5,3).Value = 0.9811433272071707
$ws.Cells.Item(6,1).Value = 0.001517954998238078
$ws.Cells.Item(6,2).Value = 0.25388421980753506
$ws.Cells.Item(6,3).Value = -0.9672334251592749
$ws.Cells.Item(7,1).Value = -0.001903814559868042
$ws.Cells.Item(7,2).Value = -0.31842254037326667
$ws.Cells.Item(7,3).Value = 0.9479469717617948
$ws.Cells.Item(8,1).Value = 0.002245110873152343
$ws.Cells.Item(8,2).Value = 0.3755039216889036
$ws.Cells.Item(8,3).Value = -0.9268180858579643
$ws.Cells.Item(9,1).Value = -0.0026193129364291043
$ws.Cells.Item(9,2).Value = -0.43809296220881694
$ws.Cells.Item(9,3).Value = 0.8989258565993332
$ws.Cells.Item(10,1).Value = 0.002934810274066609
$ws.Cells.Item(10,2).Value = 0.490858766041759
$ws.Cells.Item(10,3).Value = -0.8712342157471875
$ws.Cells.Item(11,1).Value = -0.0032900104100215817
$ws.Cells.Item(11,2).Value = -0.5502699842317609
$ws.Cells.Item(11,3).Value = 0.8349803113158295
$ws.Cells.Item(12,1).Value = 0.003575744673273488
$ws.Cells.Item(12,2).Value = 0.598057355371699
$ws.Cells.Item(12,3).Value = -0.8014453279768001
$ws.Cells.Item(13,1).Value = -0.003904786289461346
$ws.Cells.Item(13,2).Value = -0.6530936127410855
$ws.Cells.Item(13,3).Value = 0.7572671164396554
$ws.Cells.Item(14,1).Value = 0.004157541817403305
$ws.Cells.Item(14,2).Value = 0.6953649111688653
$ws.Cells.Item(14,3).Value = -0.7186448045878802
$ws.Cells.Item(15,1).Value = -0.004453662147187283
$ws.Cells.Item(15,2).Value = -0.7448949370092951
$ws.Cells.Item(15,3).Value = 0.6671669189276377
$ws.Cells.Item(16,1).Value = 0.004670919699352832
$ws.Cells.Item(16,2).Value = 0.7812290238136452
$ws.Cells.Item(16,3).Value = -0.6242270379119613
$ws.Cells.Item(17,1).Value = -0.004927942843629012
$ws.Cells.Item(17,2).Value = -0.824219712269636
$ws.Cells.Item(17,3).Value = 0.5662486920828059
$ws.Cells.Item(18,1).Value = 0.005107801621124121
$ws.Cells.Item(18,2).Value = 0.8542988912838118
$ws.Cells.Item(18,3).Value = -0.5197569765898763
$ws.Cells.Item(19,1).Value = -0.005320311079320032
$ws.Cells.Item(19,2).Value = -0.8898441660352499
$ws.Cells.Item(19,3).Value = 0.45623355253975995
$ws.Cells.Item(20,1).Value = 0.005461395919315531
$ws.Cells.Item(20,2).Value = 0.9134386449448185
$ws.Cells.Item(20,3).Value = -0.4069398175111232
$ws.Cells.Item(21,1).Value = -0.005624881440795569
$ws.Cells.Item(21,2).Value = -0.9407840268366269
$ws.Cells.Item(21,3).Value = 0.33895984357713776
$ws.Cells.Item(22,1).Value = 0.005726246888487177
$ws.Cells.Item(22,2).Value = 0.9577358567432487
$ws.Cells.Item(22,3).Value = -0.2875921396784821
$ws.Cells.Item(23,1).Value = -0.005837222709232087
$ws.Cells.Item(23,2).Value = -0.9762982435890242
$ws.Cells.Item(23,3).Value = 0.21635079476634533
$ws.Cells.Item(24,1).Value = 0.005898265112868237
$ws.Cells.Item(24,2).Value = 0.9865066025052299
$ws.Cells.Item(24,3).Value = -0.16361520003424695
$ws.Cells.Item(25,1).Value = -0.005954358401337208
$ws.Cells.Item(25,2).Value = -0.9958890646860539
$ws.Cells.Item(25,3).Value = 0.09038537743885801
$ws.Cells.Item(26,1).Value = 0.0059747456582586855
$ws.Cells.Item(26,2).Value = 0.9992984963957261
$ws.Cells.Item(26,3).Value = -0.036970495203068214
$ws.Cells.Item(27,1).Value = -0.005974755110647236
$ws.Cells.Item(27,2).Value = -0.9993000786464137
$ws.Cells.Item(27,3).Value = -0.03692770123686471
$ws.Cells.Item(28,1).Value = 0.005954381438418042
$ws.Cells.Item(28,2).Value = 0.9958929211186931
$ws.Cells.Item(28,3).Value = 0.09034287469060583
$ws.Cells.Item(29,1).Value = -0.005898306426805656
$ws.Cells.Item(29,2).Value = -0.9865135183259799
$ws.Cells.Item(29,3).Value = -0.1635734946175323
$ws.Cells.Item(30,1).Value = 0.0058372767756569915
$ws.Cells.Item(30,2).Value = 0.9763072942847668
$ws.Cells.Item(30,3).Value = 0.21630994735841966
$ws.Cells.Item(31,1).Value = -0.005726317393330945
$ws.Cells.Item(31,2).Value = -0.957747659074566
$ws.Cells.Item(31,3).Value = -0.28755283136547244
$ws.Cells.Item(32,1).Value = 0.005624963070010575
$ws.Cells.Item(32,2).Value = 0.9407976915035984
$ws.Cells.Item(32,3).Value = 0.33892191350215356
$ws.Cells.Item(33,1).Value = -0.005461491070674445
$ws.Cells.Item(33,2).Value = -0.913454573066334
$ws.Cells.Item(33,3).Value = -0.4069040612472263
$ws.Cells.Item(34,1).Value = 0.005320415004345975
$ws.Cells.Item(34,2).Value = 0.8898615629933536
$ws.Cells.Item(34,3).Value = 0.45619961846894086
$ws.Cells.Item(35,1).Value = -0.00510791528450765
$ws.Cells.Item(35,2).Value = -0.8543179182957685
$ws.Cells.Item(35,3).Value = -0.5197257004230508
$ws.Cells.Item(36,1).Value = 0.004928062318662813
$ws.Cells.Item(36,2).Value = 0.8242397122846975
$ws.Cells.Item(36,3).Value = 0.5662195783392008
$ws.Cells.Item(37,1).Value = -0.0046710445338058555
$ws.Cells.Item(37,2).Value = -0.7812499208555906
$ws.Cells.Item(37,3).Value = -0.6242008831346657
$ws.Cells.Item(38,1).Value = 0.004453789369298145
$ws.Cells.Item(38,2).Value = 0.7449162338831113
$ws.Cells.Item(38,3).Value = 0.667143139256978
$ws.Cells.Item(39,1).Value = -0.004157669738647154
$ws.Cells.Item(39,2).Value = -0.6953863249476983
$ws.Cells.Item(39,3).Value = -0.7186240831325362
$ws.Cells.Item(40,1).Value = 0.0039049128978508986
$ws.Cells.Item(40,2).Value = 0.653114806886069
$ws.Cells.Item(40,3).Value = 0.757248836698633
$ws.Cells.Item(41,1).Value = -0.0035758673702375792
$ws.Cells.Item(41,2).Value = -0.5980778946276303
$ws.Cells.Item(41,3).Value = -0.8014300001437002
$ws.Cells.Item(42,1).Value = 0.0032901280334741217
$ws.Cells.Item(42,2).Value = 0.5502896743132761
$ws.Cells.Item(42,3).Value = 0.8349673343321353
$ws.Cells.Item(43,1).Value = -0.0029349197477157862
$ws.Cells.Item(43,2).Value = -0.49087709174710326
$ws.Cells.Item(43,3).Value = -0.8712238903083296
$ws.Cells.Item(44,1).Value = 0.002619413753734472
$ws.Cells.Item(44,2).Value = 0.43810983895991
$ws.Cells.Item(44,3).Value = 0.8989176311977135
$ws.Cells.Item(45,1).Value = -0.0022451999629197385
$ws.Cells.Item(45,2).Value = -0.37551883517264256
$ws.Cells.Item(45,3).Value = -0.9268120432470157
$ws.Cells.Item(46,1).Value = 0.0019038918341658597
$ws.Cells.Item(46,2).Value = 0.31843547604452427
$ws.Cells.Item(46,3).Value = 0.9479426263187983
$ws.Cells.Item(47,1).Value = -0.001518017859968651
$ws.Cells.Item(47,2).Value = -0.25389474276327534
$ws.Cells.Item(47,3).Value = -0.9672306628818934
$ws.Cells.Item(48,1).Value = 0.001155637470468096
$ws.Cells.Item(48,2).Value = 0.19328625820424777
$ws.Cells.Item(48,3).Value = 0.981141726200368
$ws.Cells.Item(49,1).Value = -0.000765427395482992
$ws.Cells.Item(49,2).Value = -0.12802091024066778
$ws.Cells.Item(49,3).Value = -0.9917711735385604
$ws.Cells.Item(50,1).Value = 0.0003874447484656346
$ws.Cells.Item(50,2).Value = 0.06480212386131848
$ws.Cells.Item(50,3).Value = 0.9978980582352233
$ws.Cells.Item(51,1).Value = 0
$ws.Cells.Item(51,2).Value = 0
$ws.Cells.Item(51,3).Value = -1
